# iecorrect: finished testing v1!
#
# - "string" and "numeric" sheets: rename the first column header from
#   "varname" to "id", and the old "idvalue" column becomes "varname".
# - "drop" sheet: rename "idvalue" header to "id" and insert a new
#   "n_obs" column right after it (pushing "initials"/"notes" over).

$wb = $excel.ActiveWorkbook

$sString = $wb.Worksheets.Item("string")
$sString.Range("A1").Value = "id"
$sString.Range("B1").Value = "varname"

$sNumeric = $wb.Worksheets.Item("numeric")
$sNumeric.Range("A1").Value = "id"
$sNumeric.Range("B1").Value = "varname"

$sDrop = $wb.Worksheets.Item("drop")
$sDrop.Range("A1").Value = "id"
$sDrop.Columns.Item(2).Insert()
$sDrop.Range("B1").Value = "n_obs"
